$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.170.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "'1.878.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'314.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").Value = "'0.5134"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("D8").Value = "'0.3903"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.31%  "

$ws.Range("D9").Value = "'0.08351"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("D10").Value = "'1.120"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.87%  "

$ws.Range("D11").Value = "'41.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.15%  "

$ws.Range("D12").Value = "'6.230"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D13").Value = "'20.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").Value = "'1.876.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").Value = "'7.259"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.06%  "

$ws.Range("D16").Value = "'1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D18").Value = "'91.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").Value = "'0.06669"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("D20").Value = "'17.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.21%  "

$ws.Range("D21").Value = "'1.005"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D22").Value = "'6.046"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("D23").Value = "'28.211.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.66%  "

$ws.Range("D24").Value = "'11.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("D25").Value = "'2.266"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "

$ws.Range("D26").Value = "'2.091.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'159.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.25%  "

$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.483"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.16%  "

$ws.Range("D29").Value = "'20.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.09%  "

$ws.Range("D30").Value = "'125.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").Value = "'0.1062"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "

$ws.Range("D32").Value = "'1.038"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("D33").Value = "'5.833"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.19%  "

$ws.Range("D34").Value = "'3.613"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.22%  "

$ws.Range("D35").Value = "'9.632"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").Value = "'0.02448"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.33%  "

$ws.Range("D37").Value = "'0.06564"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("D38").Value = "'0.2188"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("E39").Value = "  -0.65%  "

$ws.Range("D40").Value = "'0.6497"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.30%  "

$ws.Range("D41").Value = "'4.994"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.38%  "

$ws.Range("D42").Value = "'1.227"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.26%  "

$ws.Range("D43").Value = "'11.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.67%  "

$ws.Range("D44").Value = "'0.6136"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").Value = "'13.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("D46").Value = "'1.285"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").Value = "'3.675"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("D48").Value = "'2.018"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("D49").Value = "'1.229"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.08%  "

$ws.Range("D50").Value = "'120.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").Value = "'0.06905"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.76%  "
